$d = $word.ActiveDocument

function Replace-InRange($rng, $old, $new) {
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1) Ativação date 2015 -> 2025
Replace-InRange $d.Content "Ativação: 01/01/2015" "Ativação: 01/01/2025"

# 2) Remove the "8426375 - Wendell de Queiróz Lamas" teacher line (and the
#    manual line break that preceded it), leaving only Ana Karine's entry.
Replace-InRange $d.Content "^l8426375 - Wendell de Queiróz Lamas" ""

# 3) "Programa resumido" paragraph (Portuguese): join the two runs (drop the
#    manual line break) and trim the "Energia solar. Energia eólica." /
#    "Energia da biomassa." mentions.
#    NOTE: Find.Execute collapses the range to the replaced text, so each
#    paragraph's Range must be re-fetched from the Paragraph object between
#    calls rather than reusing the same Range variable.
$pResumoPtPara = $d.Paragraphs.Item(11)
Replace-InRange $pResumoPtPara.Range "antrópicas,^ldemanda" "antrópicas,demanda"
Replace-InRange $pResumoPtPara.Range `
    "nucleares. Energia solar. Energia eólica. Energia fóssil. Energia da biomassa. Impactos" `
    "nucleares. Energia fóssil. Impactos"

# 4) "Programa resumido" paragraph (English, italic): trim "Solar energy.
#    Wind energy." / "Biomass energy." mentions.
$pResumoEnPara = $d.Paragraphs.Item(12)
Replace-InRange $pResumoEnPara.Range `
    "Power plants. Solar energy. Wind energy. Fossil energy. Biomass energy. Environmental" `
    "Power plants. Fossil energy. Environmental"

# 5) "Programa" paragraph (Portuguese): join the two runs, trim the same
#    "solar/eólica/biomassa" mentions, and append the new sentence about
#    didactic trips (note the source text's "regiona." truncation is kept
#    verbatim as it appears in the target).
$pProgPtPara = $d.Paragraphs.Item(14)
Replace-InRange $pProgPtPara.Range "antrópicas,^ldemanda" "antrópicas,demanda"
Replace-InRange $pProgPtPara.Range `
    "nucleares. Energia solar. Energia eólica. Energia fóssil. Energia da biomassa. Impactos ambientais decorrentes da geração, transmissão, disponibilidade e oferta de energia no desenvolvimento regional." `
    "nucleares. Energia fóssil. Impactos ambientais decorrentes da geração, transmissão, disponibilidade e oferta de energia no desenvolvimento regiona. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina."

# 6) "Programa" paragraph (English, italic): trim "Solar energy. Wind
#    energy." / "Biomass energy." mentions and append the new sentence
#    about didactic trips.
$pProgEnPara = $d.Paragraphs.Item(15)
Replace-InRange $pProgEnPara.Range `
    "nuclear. Solar energy. Wind energy. Fossil energy. Biomass energy. Environmental impacts of energy generation, transmission, availability, and supply in regional development." `
    "nuclear. Fossil energy. Environmental impacts of energy generation, transmission, availability, and supply in regional development. The discipline may have didactic trips to complement the content of the discipline."

# 7) Método / Critério / Norma de recuperação text replacements.
Replace-InRange $d.Content `
    "Os alunos efetuarão monografias em grupos a serem selecionados em classe." `
    "O método de avaliação será composto por avaliação teórica, apresentação escrita e oral."

Replace-InRange $d.Content `
    "Dois seminários, pesos 1 e 2." `
    "Para o cálculo da nota final (NF) será adotada a média ponderada de provas e atividades."

Replace-InRange $d.Content `
    "Para os alunos reprovados por nota, mas beneficiados pelo sistema de recuperação, esta será realizada através da aplicação de uma única prova teórica, abrangendo todo o programa do semestre letivo." `
    "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 => 5,0 Aprovado"
